# Update the marksheet's correct/total mark figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: number of right answers used for marking went from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row: total score went from 21 to 35
$ws.Range("B12").Value = 35

# Corr/total marks text e.g. "21/84" -> "35/140"
$ws.Range("E12").Value = "35/140"
